$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab from "Through 2022-10-11" to "Through 2022-10-12"
$ws.Name = "Through 2022-10-12"

# Update the header label in I1 (shared string "2022 (through 10-11)" -> "2022 (through 10-12)")
$ws.Range("I1").Value = "2022 (through 10-12)"

# Update November's value for 2022 (I11): 37 -> 39
$ws.Range("I11").Value = 39

# Update the Total for 2022 (I14): 1315 -> 1317
$ws.Range("I14").Value = 1317
